$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" A1 conversion note text ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 1.7 = 6217.42 pesos`n✅ 6217.42 pesos = 1.69 = 933.2 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update "tasas" sheet rate figures ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 588.99
$ws2.Range("O10").Value = 3662
$ws2.Range("N12").Value = 3671
$ws2.Range("O12").Value = 550.999
